# Apply "Corr/total marks" changes to the marksheet (quiz sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking (points per correct answer): 3 -> 5
$ws.Range("B11").Value = 5

# Total marks obtained: 81 -> 135 (27 correct * 5 marks)
$ws.Range("B12").Value = 135

# Total / Max marks text: "80/84" -> "135/140"
$ws.Range("E12").Value = "135/140"
